# Helper: classic VBA-style RGB() -> packed BGR-order long used by Excel's
# Font.Color / Interior.Color (same convention as real Excel COM).
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet to "Sample Data" and insert a brand-new
#    "Instructions" sheet in front of it (so tab order is Instructions, then
#    Sample Data). NOTE: the handle returned by Worksheets.Add() can alias
#    the "before" sheet handle across the call, so re-fetch both sheets by
#    name immediately afterwards instead of trusting pre-Add references.
# ---------------------------------------------------------------------------
$wb.ActiveSheet.Name = "Sample Data"
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Sample Data"))
$newSheet.Name = "Instructions"

$dataSheet = $wb.Worksheets.Item("Sample Data")
$instrSheet = $wb.Worksheets.Item("Instructions")

# ===========================================================================
# 2. "Sample Data" sheet — update headers, trim sample rows, restyle.
# ===========================================================================

# Drop the last two demo rows (Beach Sunset / Laptop Computer) -> 9 rows left.
$dataSheet.Rows.Item(11).Delete()
$dataSheet.Rows.Item(10).Delete()

# Rename the URL / alt-text headers to the machine-friendly column names the
# tool looks for.
$dataSheet.Range("C1").Value = "image_url"
$dataSheet.Range("D1").Value = "alt_text"

# Header row: bold white text on a blue fill.
$headerRange = $dataSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Font.Color = RGB(255, 255, 255)
$headerRange.Interior.Color = RGB(68, 114, 196)

# Alt-text column: green if an example value is already filled in, otherwise
# yellow to flag that it still needs to be generated.
foreach ($r in 2..9) {
    $cell = $dataSheet.Range("D$r")
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) {
        $cell.Interior.Color = RGB(255, 242, 204)
    } else {
        $cell.Interior.Color = RGB(232, 245, 232)
    }
}

# Re-size columns: URL / alt-text columns a bit wider, price column narrower.
$dataSheet.Columns.Item(3).ColumnWidth = 60
$dataSheet.Columns.Item(4).ColumnWidth = 60
$dataSheet.Columns.Item(6).ColumnWidth = 9

# ===========================================================================
# 3. "Instructions" sheet — write the how-to-use copy.
# ===========================================================================

$instrSheet.Columns.Item(1).ColumnWidth = 59
$instrSheet.Columns.Item(2).ColumnWidth = 2
$instrSheet.Columns.Item(3).ColumnWidth = 2
$instrSheet.Columns.Item(4).ColumnWidth = 2

$instrSheet.Range("A1").Value = "SEO Alt Text Generator - Instructions"

$instrSheet.Range("A3").Value = "How to use this tool:"

$instrSheet.Range("A4").Value = "1. Your Excel file must have a column with image URLs"
$instrSheet.Range("A5").Value = "2. Column names that work for image URLs:"
$instrSheet.Range("A6").Value = "   - image_url"
$instrSheet.Range("A7").Value = "   - image"
$instrSheet.Range("A8").Value = "   - url"
$instrSheet.Range("A9").Value = "   - image_link"

$instrSheet.Range("A11").Value = "3. Optional: Alt text column (will be created if missing)"
$instrSheet.Range("A12").Value = "   Column names that work for alt text:"
$instrSheet.Range("A13").Value = "   - alt_text"
$instrSheet.Range("A14").Value = "   - alt"
$instrSheet.Range("A15").Value = "   - description"
$instrSheet.Range("A16").Value = "   - alt_description"

$instrSheet.Range("A18").Value = "4. You can include any other columns with additional data"

$instrSheet.Range("A20").Value = "5. The tool will:"
$instrSheet.Range("A21").Value = "   - Automatically detect your image URL column"
$instrSheet.Range("A22").Value = "   - Create an alt_text column if one doesn't exist"
$instrSheet.Range("A23").Value = "   - Generate SEO-friendly alt text for missing entries"
$instrSheet.Range("A24").Value = "   - Allow you to edit any generated text"
$instrSheet.Range("A25").Value = "   - Export the updated file for download"

$instrSheet.Range("A27").Value = "See the 'Sample Data' sheet for an example format"

# Title banner: big bold white text on blue fill, spanning A:D.
$titleRange = $instrSheet.Range("A1:D1")
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 16
$titleRange.Font.Color = RGB(255, 255, 255)
$titleRange.Interior.Color = RGB(68, 114, 196)

# Numbered step headings: bold, size 12, light-green fill.
foreach ($r in @(4, 5, 11, 18, 20)) {
    $cell = $instrSheet.Range("A$r")
    $cell.Font.Bold = $true
    $cell.Font.Size = 12
    $cell.Interior.Color = RGB(232, 245, 232)
}

# Bulleted sub-items: italic, light-gray fill.
foreach ($r in @(6, 7, 8, 9, 13, 14, 15, 16, 21, 22, 23, 24, 25)) {
    $cell = $instrSheet.Range("A$r")
    $cell.Font.Italic = $true
    $cell.Interior.Color = RGB(248, 248, 248)
}

[void]$instrSheet.Activate()
[void]$instrSheet.Range("A1").Select()

Write-Host "done"
